$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39
$ws.Range("A39").Value = 111880562
$ws.Range("B39").Value = 90792
$ws.Range("E39").Value = 4361
$ws.Range("F39").Value = "Orange taggsvamp"
$ws.Range("G39").Value = "Hydnellum aurantiacum"
$ws.Range("H39").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = "3"
$ws.Range("I39").Style = "Normal"
$ws.Range("Q39").Value = 509658
$ws.Range("R39").Value = 6753521

# Row 40
$ws.Range("A40").Value = 111880484
$ws.Range("B40").Value = 90792
$ws.Range("E40").Value = 4361
$ws.Range("F40").Value = "Orange taggsvamp"
$ws.Range("G40").Value = "Hydnellum aurantiacum"
$ws.Range("H40").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I40").NumberFormat = "@"
$ws.Range("I40").Value = "11"
$ws.Range("I40").Style = "Normal"
$ws.Range("Q40").Value = 509901
$ws.Range("R40").Value = 6753525
$ws.Range("AJ40").Value = "tall"
$ws.Range("AK40").Value = "Pinus sylvestris"
$ws.Range("AO40").Value = "Pinus sylvestris"

# Row 41
$ws.Range("A41").Value = 111880462
$ws.Range("B41").Value = 89100
$ws.Range("I41").NumberFormat = "@"
$ws.Range("I41").Value = "1"
$ws.Range("I41").Style = "Normal"
$ws.Range("Q41").Value = 509970
$ws.Range("R41").Value = 6753250
$ws.Range("AJ41").Value = "tall"
$ws.Range("AK41").Value = "Pinus sylvestris"
$ws.Range("AL41").Value = "vid tallar"
$ws.Range("AO41").Value = "Pinus sylvestris # vid tallar"

# Row 42
$ws.Range("A42").Value = 111880574
$ws.Range("B42").Value = 90792
$ws.Range("I42").NumberFormat = "@"
$ws.Range("I42").Value = "2"
$ws.Range("I42").Style = "Normal"
$ws.Range("Q42").Value = 509596
$ws.Range("R42").Value = 6753392

# Row 43
$ws.Range("A43").Value = 111880500
$ws.Range("B43").Value = 89100
$ws.Range("E43").Value = 5754
$ws.Range("F43").Value = "Gultoppig fingersvamp"
$ws.Range("G43").Value = "Ramaria testaceoflava"
$ws.Range("H43").Value = "(Bres.) Corner"
$ws.Range("I43").NumberFormat = "@"
$ws.Range("I43").Value = "4"
$ws.Range("I43").Style = "Normal"
$ws.Range("Q43").Value = 509899
$ws.Range("R43").Value = 6753571
$ws.Range("AJ43").Value = "gran"
$ws.Range("AK43").Value = "Picea abies"
$ws.Range("AO43").Value = "Picea abies"

# Row 44
$ws.Range("B44").Value = 90786

# Row 45
$ws.Range("A45").Value = 111880509
$ws.Range("B45").Value = 90786
$ws.Range("E45").Value = 3100
$ws.Range("F45").Value = "Talltaggsvamp"
$ws.Range("G45").Value = "Bankera fuligineoalba"
$ws.Range("H45").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("I45").NumberFormat = "@"
$ws.Range("I45").Value = "6"
$ws.Range("I45").Style = "Normal"
$ws.Range("Q45").Value = 509834
$ws.Range("R45").Value = 6753644

# Row 46
$ws.Range("A46").Value = 111880580
$ws.Range("B46").Value = 90792
$ws.Range("E46").Value = 4361
$ws.Range("F46").Value = "Orange taggsvamp"
$ws.Range("G46").Value = "Hydnellum aurantiacum"
$ws.Range("H46").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I46").NumberFormat = "@"
$ws.Range("I46").Value = "3"
$ws.Range("I46").Style = "Normal"
$ws.Range("Q46").Value = 509755
$ws.Range("R46").Value = 6753236
$ws.Range("AL46").ClearContents()
$ws.Range("AO46").Value = "Pinus sylvestris"

# Row 47
$ws.Range("A47").Value = 111880601
$ws.Range("B47").Value = 89100
$ws.Range("E47").Value = 5754
$ws.Range("F47").Value = "Gultoppig fingersvamp"
$ws.Range("G47").Value = "Ramaria testaceoflava"
$ws.Range("H47").Value = "(Bres.) Corner"
$ws.Range("I47").NumberFormat = "@"
$ws.Range("I47").Value = "4"
$ws.Range("I47").Style = "Normal"
$ws.Range("Q47").Value = 509942
$ws.Range("R47").Value = 6753225

# Row 48
$ws.Range("A48").Value = 111880591
$ws.Range("B48").Value = 90792
$ws.Range("E48").Value = 4361
$ws.Range("F48").Value = "Orange taggsvamp"
$ws.Range("G48").Value = "Hydnellum aurantiacum"
$ws.Range("H48").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I48").NumberFormat = "@"
$ws.Range("I48").Value = "8"
$ws.Range("I48").Style = "Normal"
$ws.Range("Q48").Value = 509822
$ws.Range("R48").Value = 6753234

# Row 49
$ws.Range("A49").Value = 111880475
$ws.Range("B49").Value = 89100
$ws.Range("E49").Value = 5754
$ws.Range("F49").Value = "Gultoppig fingersvamp"
$ws.Range("G49").Value = "Ramaria testaceoflava"
$ws.Range("H49").Value = "(Bres.) Corner"
$ws.Range("I49").NumberFormat = "@"
$ws.Range("I49").Value = "2"
$ws.Range("I49").Style = "Normal"
$ws.Range("Q49").Value = 509958
$ws.Range("R49").Value = 6753363
$ws.Range("AJ49").Value = "gran"
$ws.Range("AK49").Value = "Picea abies"
$ws.Range("AO49").Value = "Picea abies"
